$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A6: copy format from A1, set value 6
$ws.Range("A1").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 6

# B6: copy format from B1 ("JuanPa" style), set value "Sutano"
$ws.Range("B1").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = "Sutano"

# C6: copy format from C1 ("Host" style), set value "Host" (reuse shared string)
$ws.Range("C1").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Host"

# D6: copy format from A1 (general numfmt), then set horizontal alignment left, value 1
$ws.Range("A1").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").HorizontalAlignment = -4131
$ws.Range("D6").Value = 1

Write-Host "done"
